# Add files via upload
# Duplicate the existing distance table (rows 2:37) into rows 38:73, swapping
# the "from" / "to" columns (B <-> C), and touch column E on row 1 so the
# sheet's used range grows to A1:E73 the same way Excel would record it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- give the new block (A38:D73) the same cell formatting as the source
# block (A2:D37) before filling in values, so the duplicated rows end up
# styled (s="1") just like the originals.
$ws.Range("A2:D37").Copy()
$ws.Range("A38:D73").PasteSpecial(-4122)
$ws.Range("A38:D73").RowHeight = 15.75

# --- copy the 36 data rows down to rows 38-73 with the B/C (state pair)
# columns swapped, keeping the same weight (A) and distance (D).
for ($i = 2; $i -le 37; $i++) {
    $weight = $ws.Cells.Item($i, 1).Value2
    $fromState = $ws.Cells.Item($i, 2).Value2
    $toState = $ws.Cells.Item($i, 3).Value2
    $distance = $ws.Cells.Item($i, 4).Value2

    $destRow = $i + 36
    $ws.Cells.Item($destRow, 1).Value2 = $weight
    $ws.Cells.Item($destRow, 2).Value2 = $toState
    $ws.Cells.Item($destRow, 3).Value2 = $fromState
    $ws.Cells.Item($destRow, 4).Value2 = $distance
}

# --- touch column E on the header row (format-only, no value) so the used
# range / dimension extends to column E, matching the saved workbook.
$ws.Range("D1").Copy()
$ws.Range("E1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- restore the view: scroll near the bottom of the data and select I65.
$ws.Range("I65").Select() | Out-Null

